$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the field headers for consistency: R4_Month -> Usage_Date, R4_Count -> Usage_Count
$ws.Range("K1").Value = "Usage_Date"
$ws.Range("L1").Value = "Usage_Count"

# Update the selection to match the edited header cells
$ws.Range("K1:L1").Select()
